$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New translation rows to append (label key -> translation text), mirrors the
# "cs" column A used throughout the sheet.
$rows = @(
    @("lab.atomizer.purchase.button", "Pořídit atomizér"),
    @("lab.atomizer.purchase.submit", "Pořídit"),
    @("lab.atomizer.purchase.driptipId.label", "Náústek"),
    @("lab.atomizer.purchase.driptipId.label.tooltip", "Můžete si vybrat náústek k atomizéru, který se pak bude nabízet jako výchozí ve zbytku aplikace (např. můžete vytvořit a přiřadit výchozí náústek výrobce - něco jako Taifun GT IV Native nebo tak)."),
    @("lab.atomizer.purchase.success", "Úspěšně jste si pořídili atomizér [{{data.atomizer.name}}]."),
    @("lab.atomizer.user.button.delete", "Odstranit atomizér"),
    @("lab.atomizer.user.button.delete.confirm.title", "Odstranit atomizér"),
    @("lab.atomizer.user.button.delete.confirm.ok", "Odstranit atomizér"),
    @("lab.atomizer.user.button.delete.confirm", "Opravdu si přejete odstranit vybraný atomizér? Tato akce pouze odstraní záznam o jeho vlastnictví a přidružená data (např. výchozí náústek). Dále se přestane nabízet v různých seznamech. Tato akce je nezvratná, nikoli však destruktivní."),
    @("lab.atomizer.user.deleted.success", "Atomizér [{{data.atomizer.name}}] byl úspěšně odstraněn."),
    @("lab.atomizer.user.edit.button", "Upravit atomizér"),
    @("lab.atomizer.user.driptipId.label", "Náústek"),
    @("lab.atomizer.user.update", "Uložit"),
    @("lab.atomizer.user.update.success", "Atomizér byl úspěšně uložen.")
)

$startRow = 1148
$lastRow = $startRow + $rows.Count - 1

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $srcRow = $r - 1

    # Copy the row above so the new row inherits the same styling (style id 1
    # on columns A:C, matching the rest of the translation table).
    $ws.Range("A" + $srcRow + ":C" + $srcRow).Copy($ws.Range("A" + $r + ":C" + $r))

    $ws.Cells.Item($r, 1).Value = "cs"
    $ws.Cells.Item($r, 2).Value = $rows[$i][0]
    $ws.Cells.Item($r, 3).Value = $rows[$i][1]
}

# Match the row heights for the two long, wrapped tooltip/confirm strings.
$ws.Rows.Item(1151).RowHeight = 39
$ws.Rows.Item(1156).RowHeight = 39

# Leave the view scrolled/selected near the newly-added rows, like the author did.
[void]$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1129
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Cells.Item(1154, 2).Select()
